$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple +4 increments to column B for these rows
$simpleRows = @(8, 9, 10, 11, 12, 13, 17, 18, 19, 20, 21, 22)
foreach ($r in $simpleRows) {
    $cell = $ws.Cells.Item($r, 2)  # column B
    $old = $cell.Value()
    $cell.Value = $old + 4
}

# Rows 14, 15, 16 get a 3-way rotation of columns A,D,E,F,G,H,Q,R plus B+4
# Capture original (before) values first
function Get-RowData($r) {
    return @{
        A = $ws.Cells.Item($r, 1).Value()
        B = $ws.Cells.Item($r, 2).Value()
        D = $ws.Cells.Item($r, 4).Value()
        E = $ws.Cells.Item($r, 5).Value()
        F = $ws.Cells.Item($r, 6).Value()
        G = $ws.Cells.Item($r, 7).Value()
        H = $ws.Cells.Item($r, 8).Value()
        Q = $ws.Cells.Item($r, 17).Value()
        R = $ws.Cells.Item($r, 18).Value()
    }
}

function Set-RowData($r, $data, $bBonus) {
    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B + $bBonus
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 17).Value = $data.Q
    $ws.Cells.Item($r, 18).Value = $data.R
}

$row14 = Get-RowData 14
$row15 = Get-RowData 15
$row16 = Get-RowData 16

# New row14 = old row15 (B+4)
Set-RowData 14 $row15 4
# New row15 = old row16 (B+4)
Set-RowData 15 $row16 4
# New row16 = old row14 (B+4)
Set-RowData 16 $row14 4
